$wb = $excel.ActiveWorkbook

# Rename the single demo worksheet from the old, verbose auto-generated
# name to a short, reusable one ("line-items") now that the module
# supports exporting/importing multiple sheets per workbook.
$ws = $wb.Worksheets.Item(1)
$ws.Name = "line-items"
